# "se cambia data para ruta critica y cambios en el feature"
#
# Data-driven test fixture edits for ConsultaMovimientos.xlsx:
#   - TarjetasCredito: "usuario" column value updated (pagotdc1 -> zutarjeta7)
#     and the second credit-card row's franquicia/mascara updated
#     (Personal Visa / *8298 -> Personal American Express / *9105).
#   - Eprepago: "usuario" column value updated (autotest28 -> invictus10)
#     and the card mask updated (*5214 -> *9344).
# Final view state left pointing at the Eprepago tab (cell A2 selected),
# matching where the author was working when the fixture was saved.

$wb = $excel.ActiveWorkbook

# ---- TarjetasCredito --------------------------------------------------
$wsTarjetas = $wb.Worksheets.Item("TarjetasCredito")
$wsTarjetas.Range("D2").Value = "zutarjeta7"
$wsTarjetas.Range("D3").Value = "zutarjeta7"
$wsTarjetas.Range("M3").Value = "Personal American Express"
$wsTarjetas.Range("N3").Value = "*9105"

# ---- Eprepago ----------------------------------------------------------
$wsEprepago = $wb.Worksheets.Item("Eprepago")
$wsEprepago.Range("D2").Value = "invictus10"
$wsEprepago.Range("N2").Value = "*9344"

# ---- Restore per-sheet selections -------------------------------------
$wsDepositos = $wb.Worksheets.Item("Depositos")
[void]$wsDepositos.Range("A1").Select()

[void]$wsTarjetas.Range("M3").Select()

$wsInversiones = $wb.Worksheets.Item("Inversiones")
[void]$wsInversiones.Range("M2").Select()

# ---- Activate Eprepago last, with A2 selected --------------------------
[void]$wsEprepago.Activate()
[void]$wsEprepago.Range("A2").Select()
